$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 186; existing rows 186-285 shift down to 187-286.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new data record.
$ws.Range("A186").Value2 = 8
$ws.Range("B186").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C186").Value2 = "Coquimbo"
$ws.Range("D186").Value2 = 44529
$ws.Range("E186").Value2 = 4
$ws.Range("F186").Value2 = 100114001
$ws.Range("G186").Value2 = "Papa"
$ws.Range("H186").Value2 = "Cardinal"
$ws.Range("I186").Value2 = "1a nueva(o)"
$ws.Range("J186").Value2 = 2000
$ws.Range("K186").Value2 = 12000
$ws.Range("L186").Value2 = 13000
$ws.Range("M186").Value2 = 12500
$ws.Range("N186").Value2 = "$/saco 25 kilos"
$ws.Range("O186").Value2 = "Provincia del Elquí"
$ws.Range("P186").Value2 = 500
$ws.Range("Q186").Value2 = 25
$ws.Range("R186").Value2 = "Hortaliza"
